$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I: Pruned
$ws.Range("I1").Value = "Pruned"
$ws.Range("I2").Value = "Yes"
$ws.Range("I3").Value = "Yes"
$ws.Range("I4").Value = "Yes"
$ws.Range("I5").Value = "Yes"
$ws.Range("I6").Value = "Yes"
$ws.Range("I7").Value = "Yes"
$ws.Range("I8").Value = "Yes"

# Column J: Quadrant
$ws.Range("J1").Value = "Quadrant "
$ws.Range("J2").Value = 2
$ws.Range("J3").Value = 3
$ws.Range("J4").Value = 3
$ws.Range("J5").Value = 3
$ws.Range("J6").Value = 3
$ws.Range("J7").Value = 4
$ws.Range("J8").Value = 1

# Column K: Shade
$ws.Range("K1").Value = "Shade"
$ws.Range("K2").Value = "Dark"
$ws.Range("K3").Value = "Neutral"
$ws.Range("K4").Value = "Neutral"
$ws.Range("K5").Value = "Bright"
$ws.Range("K6").Value = "Bright"
$ws.Range("K7").Value = "Dark"
$ws.Range("K8").Value = "Neutral"

$ws.Range("L8").Select()
